$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '27.447.26'
$ws.Range("E2").Value = '  -3.42%  '
$ws.Range("D3").Value = '1.754.75'
$ws.Range("E3").Value = '  -2.79%  '
$ws.Range("E4").Value = '  +0.16%  '
$ws.Range("D5").Value = '''323.31'
$ws.Range("E5").Value = '  -1.56%  '
$ws.Range("D6").Value = '''1.002'
$ws.Range("E6").Value = '  +0.27%  '
$ws.Range("D7").Value = '''0.4412'
$ws.Range("E7").Value = '  -1.61%  '
$ws.Range("D8").Value = '''0.3704'
$ws.Range("E8").Value = '  -1.68%  '
$ws.Range("D9").Value = '''44.87'
$ws.Range("E9").Value = '  +0.44%  '
$ws.Range("D10").Value = '''0.07675'
$ws.Range("E10").Value = '  +2.07%  '
$ws.Range("D11").Value = '''1.109'
$ws.Range("E11").Value = '  -3.53%  '
$ws.Range("E12").Value = '  +0.16%  '
$ws.Range("D13").Value = '''21.50'
$ws.Range("E13").Value = '  -4.91%  '
$ws.Range("D14").Value = '''6.140'
$ws.Range("E14").Value = '  -2.65%  '
$ws.Range("D15").Value = '''7.412'
$ws.Range("E15").Value = '  -2.94%  '
$ws.Range("D16").Value = '1.756.55'
$ws.Range("E16").Value = '  -2.53%  '
$ws.Range("D17").Value = '''89.88'
$ws.Range("E17").Value = '  +11.28%  '
$ws.Range("D18").Value = '''0.00001070'
$ws.Range("E18").Value = '  -2.10%  '
$ws.Range("D19").Value = '''0.06226'
$ws.Range("E19").Value = '  -8.30%  '
$ws.Range("D20").Value = '''1.001'
$ws.Range("E20").Value = '  +0.18%  '
$ws.Range("D21").Value = '''17.31'
$ws.Range("E21").Value = '  -1.97%  '
$ws.Range("D22").Value = '''6.161'
$ws.Range("E22").Value = '  -2.68%  '
$ws.Range("D23").Value = '''0.5282'
$ws.Range("E23").Value = '  -3.00%  '
$ws.Range("D24").Value = '27.504.91'
$ws.Range("E24").Value = '  -3.27%  '
$ws.Range("D25").Value = '''11.47'
$ws.Range("E25").Value = '  -3.12%  '
$ws.Range("D26").Value = '''2.306'
$ws.Range("E26").Value = '  -4.14%  '
$ws.Range("D27").Value = '''20.42'
$ws.Range("E27").Value = '  -0.87%  '
$ws.Range("D28").Value = '''152.82'
$ws.Range("E28").Value = '  +0.66%  '
$ws.Range("D29").Value = '''2.281'
$ws.Range("E29").Value = '  -3.36%  '
$ws.Range("D30").Value = '1.956.23'
$ws.Range("E30").Value = '  -2.54%  '
$ws.Range("D31").Value = '''127.03'
$ws.Range("E31").Value = '  -4.47%  '
$ws.Range("D32").Value = '''1.170'
$ws.Range("E32").Value = '  -7.12%  '
$ws.Range("D33").Value = '''5.692'
$ws.Range("E33").Value = '  -2.43%  '
$ws.Range("D34").Value = '''0.09169'
$ws.Range("D35").Value = '''3.648'
$ws.Range("E35").Value = '  -8.90%  '
$ws.Range("D36").Value = '''12.52'
$ws.Range("E36").Value = '  +2.58%  '
$ws.Range("D37").Value = '''0.02302'
$ws.Range("E37").Value = '  -1.53%  '
$ws.Range("D38").Value = '''0.2146'
$ws.Range("E38").Value = '  -5.59%  '
$ws.Range("D39").Value = '''0.06099'
$ws.Range("E39").Value = '  -4.18%  '
$ws.Range("D40").Value = '''5.033'
$ws.Range("E40").Value = '  -2.62%  '
$ws.Range("D41").Value = '''0.6411'
$ws.Range("E41").Value = '  -2.74%  '
$ws.Range("D42").Value = '''1.175'
$ws.Range("E42").Value = '  -2.73%  '
$ws.Range("E43").Value = '  +0.26%  '
$ws.Range("D44").Value = '''7.903'
$ws.Range("E44").Value = '  -2.63%  '
$ws.Range("D45").Value = '''1.383'
$ws.Range("E45").Value = '  -4.75%  '
$ws.Range("D46").Value = '''13.62'
$ws.Range("E46").Value = '  -2.33%  '
$ws.Range("D47").Value = '''0.5937'
$ws.Range("E47").Value = '  -2.48%  '
$ws.Range("D48").Value = '''3.712'
$ws.Range("E48").Value = '  -2.64%  '
$ws.Range("D49").Value = '''125.83'
$ws.Range("E49").Value = '  -2.15%  '
$ws.Range("D50").Value = '''1.968'
$ws.Range("E50").Value = '  -3.37%  '
$ws.Range("D51").Value = '''0.06868'
$ws.Range("E51").Value = '  -3.19%  '
